$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value2 = 17.4275895051743
$ws.Cells.Item(2, 4).Value2 = 11.32734066696509
$ws.Cells.Item(2, 5).Value2 = 17.00080769389358
$ws.Cells.Item(2, 6).Value2 = 28.5233174388284
$ws.Cells.Item(2, 7).Value2 = 27.01840733239952
$ws.Cells.Item(2, 8).Value2 = 13.60205378025839
$ws.Cells.Item(2, 10).Value2 = 11.6778316495898
$ws.Cells.Item(2, 12).Value2 = 9.665761702194773
$ws.Cells.Item(2, 14).Value2 = 18.54029730951016
$ws.Cells.Item(2, 15).Value2 = 20.51747094201325
$ws.Cells.Item(3, 2).Value2 = 17.1649787209316
$ws.Cells.Item(3, 4).Value2 = 11.36427016680427
$ws.Cells.Item(3, 5).Value2 = 17.04599987502433
$ws.Cells.Item(3, 6).Value2 = 28.55345043556826
$ws.Cells.Item(3, 7).Value2 = 26.84349132782589
$ws.Cells.Item(3, 8).Value2 = 13.62599220362843
$ws.Cells.Item(3, 10).Value2 = 11.69536130478433
$ws.Cells.Item(3, 12).Value2 = 9.46871605199385
$ws.Cells.Item(3, 14).Value2 = 18.5504606241572
$ws.Cells.Item(3, 15).Value2 = 20.51796135505335
$ws.Cells.Item(4, 2).Value2 = 17.00390510434852
$ws.Cells.Item(4, 4).Value2 = 11.38823644595932
$ws.Cells.Item(4, 5).Value2 = 17.07560895068356
$ws.Cells.Item(4, 6).Value2 = 28.57944086727698
$ws.Cells.Item(4, 7).Value2 = 26.74506099223404
$ws.Cells.Item(4, 8).Value2 = 13.6432250756967
$ws.Cells.Item(4, 10).Value2 = 11.70703718188797
$ws.Cells.Item(4, 12).Value2 = 9.346050066575051
$ws.Cells.Item(4, 14).Value2 = 18.55887279423295
$ws.Cells.Item(4, 15).Value2 = 20.52338024046288
$ws.Cells.Item(5, 2).Value2 = 16.93838450939826
$ws.Cells.Item(5, 4).Value2 = 11.3983284394207
$ws.Cells.Item(5, 5).Value2 = 17.08814404619137
$ws.Cells.Item(5, 6).Value2 = 28.59191248382509
$ws.Cells.Item(5, 7).Value2 = 26.70724565593234
$ws.Cells.Item(5, 8).Value2 = 13.65088447337963
$ws.Cells.Item(5, 10).Value2 = 11.71202536339589
$ws.Cells.Item(5, 12).Value2 = 9.295699606535855
$ws.Cells.Item(5, 14).Value2 = 18.56284862657058
$ws.Cells.Item(5, 15).Value2 = 20.52687568823455
$ws.Cells.Item(6, 2).Value2 = 16.92751421565312
$ws.Cells.Item(6, 4).Value2 = 11.40002389272802
$ws.Cells.Item(6, 5).Value2 = 17.09025386103807
$ws.Cells.Item(6, 6).Value2 = 28.59409685327237
$ws.Cells.Item(6, 7).Value2 = 26.70110620556705
$ws.Cells.Item(6, 8).Value2 = 13.65219475910011
$ws.Cells.Item(6, 10).Value2 = 11.71286756875585
$ws.Cells.Item(6, 12).Value2 = 9.28731879378906
$ws.Cells.Item(6, 14).Value2 = 18.56354194391111
$ws.Cells.Item(6, 15).Value2 = 20.52753384223137
$ws.Cells.Item(7, 2).Value2 = 17.00302089216504
$ws.Cells.Item(7, 4).Value2 = 11.38837123093514
$ws.Cells.Item(7, 5).Value2 = 17.07577610213974
$ws.Cells.Item(7, 6).Value2 = 28.57960145492673
$ws.Cells.Item(7, 7).Value2 = 26.74454165551445
$ws.Cells.Item(7, 8).Value2 = 13.64332579511222
$ws.Cells.Item(7, 10).Value2 = 11.70710352140107
$ws.Cells.Item(7, 12).Value2 = 9.345372412969461
$ws.Cells.Item(7, 14).Value2 = 18.55892419350621
$ws.Cells.Item(7, 15).Value2 = 20.52342216964037
$ws.Cells.Item(8, 2).Value2 = 17.33705093696186
$ws.Cells.Item(8, 4).Value2 = 11.33980640884211
$ws.Cells.Item(8, 5).Value2 = 17.01600446346682
$ws.Cells.Item(8, 6).Value2 = 28.532151620213
$ws.Cells.Item(8, 7).Value2 = 26.9562605300328
$ws.Cells.Item(8, 8).Value2 = 13.60978144133725
$ws.Cells.Item(8, 10).Value2 = 11.68368687781129
$ws.Cells.Item(8, 12).Value2 = 9.598203666226123
$ws.Cells.Item(8, 14).Value2 = 18.54335187780221
$ws.Cells.Item(8, 15).Value2 = 20.51657825589448
$ws.Cells.Item(9, 2).Value2 = 17.99016640742034
$ws.Cells.Item(9, 4).Value2 = 11.25478114118827
$ws.Cells.Item(9, 5).Value2 = 16.91350304279526
$ws.Cells.Item(9, 6).Value2 = 28.49860848559365
$ws.Cells.Item(9, 7).Value2 = 27.44060564571211
$ws.Cells.Item(9, 8).Value2 = 13.56412946549811
$ws.Cells.Item(9, 10).Value2 = 11.64497755375389
$ws.Cells.Item(9, 12).Value2 = 10.07820279194218
$ws.Cells.Item(9, 14).Value2 = 18.52997211758415
$ws.Cells.Item(9, 15).Value2 = 20.54371878862171
$ws.Cells.Item(10, 2).Value2 = 18.46455649828667
$ws.Cells.Item(10, 4).Value2 = 11.19848627637612
$ws.Cells.Item(10, 5).Value2 = 16.84708780983672
$ws.Cells.Item(10, 6).Value2 = 28.51030379400151
$ws.Cells.Item(10, 7).Value2 = 27.83565611320465
$ws.Cells.Item(10, 8).Value2 = 13.54287810269518
$ws.Cells.Item(10, 10).Value2 = 11.62089235336695
$ws.Cells.Item(10, 12).Value2 = 10.41809276738021
$ws.Cells.Item(10, 14).Value2 = 18.53049823749163
$ws.Cells.Item(10, 15).Value2 = 20.58827670625683
$ws.Cells.Item(11, 2).Value2 = 18.67831317720285
$ws.Cells.Item(11, 4).Value2 = 11.17420588765245
$ws.Cells.Item(11, 5).Value2 = 16.81878874806796
$ws.Cells.Item(11, 6).Value2 = 28.52350810382923
$ws.Cells.Item(11, 7).Value2 = 28.02318237321708
$ws.Cells.Item(11, 8).Value2 = 13.53588003811813
$ws.Cells.Item(11, 10).Value2 = 11.61087212597995
$ws.Cells.Item(11, 12).Value2 = 10.5693406618905
$ws.Cells.Item(11, 14).Value2 = 18.53296260289112
$ws.Cells.Item(11, 15).Value2 = 20.61385206297208
$ws.Cells.Item(12, 2).Value2 = 18.7588935436062
$ws.Cells.Item(12, 4).Value2 = 11.16520175054263
$ws.Cells.Item(12, 5).Value2 = 16.80834655694299
$ws.Cells.Item(12, 6).Value2 = 28.52963955117063
$ws.Cells.Item(12, 7).Value2 = 28.09525422848051
$ws.Cells.Item(12, 8).Value2 = 13.53361372548845
$ws.Cells.Item(12, 10).Value2 = 11.60721164577952
$ws.Cells.Item(12, 12).Value2 = 10.62608340214866
$ws.Cells.Item(12, 14).Value2 = 18.53421355650161
$ws.Cells.Item(12, 15).Value2 = 20.62429484893795
$ws.Cells.Item(13, 2).Value2 = 18.74155642948687
$ws.Cells.Item(13, 4).Value2 = 11.16713249925907
$ws.Cells.Item(13, 5).Value2 = 16.81058329896128
$ws.Cells.Item(13, 6).Value2 = 28.52826877155088
$ws.Cells.Item(13, 7).Value2 = 28.07968617879859
$ws.Cells.Item(13, 8).Value2 = 13.53408475500343
$ws.Cells.Item(13, 10).Value2 = 11.60799404834154
$ws.Cells.Item(13, 12).Value2 = 10.61388719312334
$ws.Cells.Item(13, 14).Value2 = 18.53393004530544
$ws.Cells.Item(13, 15).Value2 = 20.62201218679897
$ws.Cells.Item(14, 2).Value2 = 18.68495027407613
$ws.Cells.Item(14, 4).Value2 = 11.17346130161638
$ws.Cells.Item(14, 5).Value2 = 16.81792417674481
$ws.Cells.Item(14, 6).Value2 = 28.52398988737556
$ws.Cells.Item(14, 7).Value2 = 28.02909090516587
$ws.Cells.Item(14, 8).Value2 = 13.53568589857474
$ws.Cells.Item(14, 10).Value2 = 11.61056829510192
$ws.Cells.Item(14, 12).Value2 = 10.57401983225559
$ws.Cells.Item(14, 14).Value2 = 18.53305916435218
$ws.Cells.Item(14, 15).Value2 = 20.61469603847116
$ws.Cells.Item(15, 2).Value2 = 18.65022783833053
$ws.Cells.Item(15, 4).Value2 = 11.17736264012
$ws.Cells.Item(15, 5).Value2 = 16.82245633258778
$ws.Cells.Item(15, 6).Value2 = 28.52151617987634
$ws.Cells.Item(15, 7).Value2 = 27.99823590379339
$ws.Cells.Item(15, 8).Value2 = 13.5367166079549
$ws.Cells.Item(15, 10).Value2 = 11.61216252155191
$ws.Cells.Item(15, 12).Value2 = 10.54952935825072
$ws.Cells.Item(15, 14).Value2 = 18.53256703819739
$ws.Cells.Item(15, 15).Value2 = 20.6103132310376
$ws.Cells.Item(16, 2).Value2 = 18.45053974264261
$ws.Cells.Item(16, 4).Value2 = 11.2000997276158
$ws.Cells.Item(16, 5).Value2 = 16.84897563007442
$ws.Cells.Item(16, 6).Value2 = 28.50959932303736
$ws.Cells.Item(16, 7).Value2 = 27.82355267315651
$ws.Cells.Item(16, 8).Value2 = 13.54338914725278
$ws.Cells.Item(16, 10).Value2 = 11.62156598136981
$ws.Cells.Item(16, 12).Value2 = 10.40813646434654
$ws.Cells.Item(16, 14).Value2 = 18.53038177248665
$ws.Cells.Item(16, 15).Value2 = 20.58671165339164
$ws.Cells.Item(17, 2).Value2 = 18.32746257536861
$ws.Cells.Item(17, 4).Value2 = 11.21438793089684
$ws.Cells.Item(17, 5).Value2 = 16.86573367127978
$ws.Cells.Item(17, 6).Value2 = 28.50430660360135
$ws.Cells.Item(17, 7).Value2 = 27.71834709140562
$ws.Cells.Item(17, 8).Value2 = 13.54816615261362
$ws.Cells.Item(17, 10).Value2 = 11.62757401990328
$ws.Cells.Item(17, 12).Value2 = 10.32049898686526
$ws.Cells.Item(17, 14).Value2 = 18.52960952885343
$ws.Cells.Item(17, 15).Value2 = 20.57358848959601
$ws.Cells.Item(18, 2).Value2 = 18.25648333403567
$ws.Cells.Item(18, 4).Value2 = 11.22273120721922
$ws.Cells.Item(18, 5).Value2 = 16.87555263323329
$ws.Cells.Item(18, 6).Value2 = 28.50200467026206
$ws.Cells.Item(18, 7).Value2 = 27.65857529203359
$ws.Cells.Item(18, 8).Value2 = 13.55116502341684
$ws.Cells.Item(18, 10).Value2 = 11.63111784923968
$ws.Cells.Item(18, 12).Value2 = 10.26977765851051
$ws.Cells.Item(18, 14).Value2 = 18.52937492323099
$ws.Cells.Item(18, 15).Value2 = 20.56654014840739
$ws.Cells.Item(19, 2).Value2 = 18.23242079874764
$ws.Cells.Item(19, 4).Value2 = 11.22557760145318
$ws.Cells.Item(19, 5).Value2 = 16.87890814528041
$ws.Cells.Item(19, 6).Value2 = 28.5013528231508
$ws.Cells.Item(19, 7).Value2 = 27.63846654303325
$ws.Cells.Item(19, 8).Value2 = 13.55222354604267
$ws.Cells.Item(19, 10).Value2 = 11.63233289166993
$ws.Cells.Item(19, 12).Value2 = 10.25255174447241
$ws.Cells.Item(19, 14).Value2 = 18.52933154897035
$ws.Cells.Item(19, 15).Value2 = 20.56423967305969
$ws.Cells.Item(20, 2).Value2 = 18.34058437940667
$ws.Cells.Item(20, 4).Value2 = 11.21285398713978
$ws.Cells.Item(20, 5).Value2 = 16.86393110936858
$ws.Cells.Item(20, 6).Value2 = 28.50479321186693
$ws.Cells.Item(20, 7).Value2 = 27.72947030970801
$ws.Cells.Item(20, 8).Value2 = 13.54763162789635
$ws.Cells.Item(20, 10).Value2 = 11.62692533409749
$ws.Cells.Item(20, 12).Value2 = 10.32986105313875
$ws.Cells.Item(20, 14).Value2 = 18.52967006161328
$ws.Cells.Item(20, 15).Value2 = 20.57493378346542
$ws.Cells.Item(21, 2).Value2 = 18.7015873180344
$ws.Cells.Item(21, 4).Value2 = 11.17159721934211
$ws.Cells.Item(21, 5).Value2 = 16.81576055425536
$ws.Cells.Item(21, 6).Value2 = 28.5252160221142
$ws.Cells.Item(21, 7).Value2 = 28.0439237253416
$ws.Cells.Item(21, 8).Value2 = 13.53520519217154
$ws.Cells.Item(21, 10).Value2 = 11.60980854606621
$ws.Cells.Item(21, 12).Value2 = 10.58574462746939
$ws.Cells.Item(21, 14).Value2 = 18.53330635755903
$ws.Cells.Item(21, 15).Value2 = 20.61682444209762
$ws.Cells.Item(22, 2).Value2 = 18.93536752237105
$ws.Cells.Item(22, 4).Value2 = 11.14574253965604
$ws.Cells.Item(22, 5).Value2 = 16.78587516685245
$ws.Cells.Item(22, 6).Value2 = 28.54515547617995
$ws.Cells.Item(22, 7).Value2 = 28.25558287057135
$ws.Cells.Item(22, 8).Value2 = 13.52932013801643
$ws.Cells.Item(22, 10).Value2 = 11.59940225977198
$ws.Cells.Item(22, 12).Value2 = 10.74986147652188
$ws.Cells.Item(22, 14).Value2 = 18.53753389660671
$ws.Cells.Item(22, 15).Value2 = 20.64861741460547
$ws.Cells.Item(23, 2).Value2 = 18.81081444485923
$ws.Cells.Item(23, 4).Value2 = 11.15944042653338
$ws.Cells.Item(23, 5).Value2 = 16.80167981610674
$ws.Cells.Item(23, 6).Value2 = 28.53391133574815
$ws.Cells.Item(23, 7).Value2 = 28.14207583512299
$ws.Cells.Item(23, 8).Value2 = 13.53225656196407
$ws.Cells.Item(23, 10).Value2 = 11.60488509406275
$ws.Cells.Item(23, 12).Value2 = 10.66256921058368
$ws.Cells.Item(23, 14).Value2 = 18.53510896564879
$ws.Cells.Item(23, 15).Value2 = 20.63124680898535
$ws.Cells.Item(24, 2).Value2 = 18.33465269330619
$ws.Cells.Item(24, 4).Value2 = 11.21354708187795
$ws.Cells.Item(24, 5).Value2 = 16.86474547270481
$ws.Cells.Item(24, 6).Value2 = 28.50457090817532
$ws.Cells.Item(24, 7).Value2 = 27.72443927911818
$ws.Cells.Item(24, 8).Value2 = 13.54787249995843
$ws.Cells.Item(24, 10).Value2 = 11.62721832540669
$ws.Cells.Item(24, 12).Value2 = 10.32562951022359
$ws.Cells.Item(24, 14).Value2 = 18.52964204242338
$ws.Cells.Item(24, 15).Value2 = 20.57432402965914
$ws.Cells.Item(25, 2).Value2 = 17.81413568840888
$ws.Cells.Item(25, 4).Value2 = 11.276694979134
$ws.Cells.Item(25, 5).Value2 = 16.93966534177923
$ws.Cells.Item(25, 6).Value2 = 28.50130183425358
$ws.Cells.Item(25, 7).Value2 = 27.30246402908497
$ws.Cells.Item(25, 8).Value2 = 13.57432221797153
$ws.Cells.Item(25, 10).Value2 = 11.65468176295306
$ws.Cells.Item(25, 12).Value2 = 9.950376492431708
$ws.Cells.Item(25, 14).Value2 = 18.53176518126311
$ws.Cells.Item(25, 15).Value2 = 20.53204213939349
